$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Retornos")

# Shift existing data rows 2..242 down by 6 rows (to 8..248), bottom-up to avoid clobbering.
for ($r = 242; $r -ge 2; $r--) {
    $dest = $r + 6
    $ws.Cells.Item($dest,1).Value = $ws.Cells.Item($r,1).Value2
    $ws.Cells.Item($dest,2).Value = $ws.Cells.Item($r,2).Value2
    $ws.Cells.Item($dest,3).Value = $ws.Cells.Item($r,3).Value2
}

# Copy the date-column style (s="2") down onto the newly-created tail rows (243..248).
$ws.Range("A242").Copy() | Out-Null
$ws.Range("A243:A248").PasteSpecial(-4122) | Out-Null

# Write the 6 new most-recent rows at the top (rows 2..7).
$ws.Cells.Item(2,1).Value = 45422
$ws.Cells.Item(2,2).Value = 0.01972386587771213
$ws.Cells.Item(2,3).Value = -0.004587012824913361
$ws.Cells.Item(3,1).Value = 45421
$ws.Cells.Item(3,2).Value = 0.01934235976789167
$ws.Cells.Item(3,3).Value = -0.009986021115067079
$ws.Cells.Item(4,1).Value = 45420
$ws.Cells.Item(4,2).Value = -0.007210626185958358
$ws.Cells.Item(4,3).Value = 0.002097360885380484
$ws.Cells.Item(5,1).Value = 45419
$ws.Cells.Item(5,2).Value = -0.002675840978593302
$ws.Cells.Item(5,3).Value = 0.0057914156274812
$ws.Cells.Item(6,1).Value = 45418
$ws.Cells.Item(6,2).Value = 0.1399003449597547
$ws.Cells.Item(6,3).Value = -0.0003346069146907826
$ws.Cells.Item(7,1).Value = 45415
$ws.Cells.Item(7,2).Value = -0.005211835911230533
$ws.Cells.Item(7,3).Value = 0.01091077862211098

# Update the Beta sheet value.
$wsBeta = $wb.Worksheets.Item("Beta")
$wsBeta.Cells.Item(2,2).Value = -0.3182495238079706
